$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 13).ClearContents()

$ws.Cells.Item(19, 8).Value = 13135
$ws.Cells.Item(19, 9).Value = 17180
$ws.Cells.Item(19, 10).Value = 1000
$ws.Cells.Item(19, 11).Value = 17180
$ws.Cells.Item(19, 12).Value = 1000
$ws.Cells.Item(19, 13).Value = -17005
$ws.Cells.Item(19, 14).Value = -1350

$ws.Cells.Item(33, 8).Value = 4546583
$ws.Cells.Item(33, 9).Value = 1335.5
$ws.Cells.Item(33, 11).Value = 1335.5
$ws.Cells.Item(33, 13).Value = -1106.5

$ws.Cells.Item(62, 8).Value = 1784.1666
$ws.Cells.Item(62, 9).Value = 1784.1666
$ws.Cells.Item(62, 11).Value = 1784.1666
$ws.Cells.Item(62, 13).Value = -1160.1666

$ws.Cells.Item(65, 8).Value = 1784.1666
$ws.Cells.Item(65, 9).Value = 1784.1666
$ws.Cells.Item(65, 11).Value = 8920.833000000001
$ws.Cells.Item(65, 13).Value = -5800.833000000001

$ws.Cells.Item(98, 8).Value = 1215.8948
$ws.Cells.Item(98, 9).Value = 1273.4667
$ws.Cells.Item(98, 11).Value = 1273.4667
$ws.Cells.Item(98, 13).Value = 224.5333000000001

$ws.Cells.Item(113, 8).Value = 4050
$ws.Cells.Item(113, 9).Value = 3100
$ws.Cells.Item(113, 10).Value = 5000
$ws.Cells.Item(113, 11).Value = 3100
$ws.Cells.Item(113, 12).Value = 5000
$ws.Cells.Item(113, 13).Value = 154
$ws.Cells.Item(113, 14).Value = -11508

$ws.Cells.Item(122, 8).Value = 1215.8948
$ws.Cells.Item(122, 9).Value = 1273.4667
$ws.Cells.Item(122, 11).Value = 3820.4001
$ws.Cells.Item(122, 13).Value = -1370.4001

$ws.Cells.Item(137, 8).Value = 1436.6364
$ws.Cells.Item(137, 9).Value = 1311.5
$ws.Cells.Item(137, 10).Value = 1999.75
$ws.Cells.Item(137, 11).Value = 3934.5
$ws.Cells.Item(137, 12).Value = 5999.25
$ws.Cells.Item(137, 13).Value = -1384.5
$ws.Cells.Item(137, 14).Value = -11099.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 14124.375
$ws.Cells.Item(45, 9).Value = 34137.332
$ws.Cells.Item(45, 10).Value = 2116.6
$ws.Cells.Item(45, 11).Value = 34137.332
$ws.Cells.Item(45, 12).Value = 2116.6
$ws.Cells.Item(45, 13).Value = -33760.332
$ws.Cells.Item(45, 14).Value = -2870.6

$ws.Cells.Item(122, 8).Value = 3664803.5
$ws.Cells.Item(122, 9).Value = 4275104
$ws.Cells.Item(122, 11).Value = 12825312
$ws.Cells.Item(122, 13).Value = -12822862

$ws.Cells.Item(123, 8).Value = 30414
$ws.Cells.Item(123, 10).Value = 30414
$ws.Cells.Item(123, 12).Value = 30414
$ws.Cells.Item(123, 14).Value = -40214

$ws.Cells.Item(132, 8).Value = 11619.9375
$ws.Cells.Item(132, 9).Value = 2261.75
$ws.Cells.Item(132, 10).Value = 20978.125
$ws.Cells.Item(132, 11).Value = 6785.25
$ws.Cells.Item(132, 12).Value = 62934.375
$ws.Cells.Item(132, 13).Value = -4255.25
$ws.Cells.Item(132, 14).Value = -67994.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2319.8
$ws.Cells.Item(86, 9).Value = 2319.8
$ws.Cells.Item(86, 11).Value = 2319.8
$ws.Cells.Item(86, 13).Value = -1196.8

$ws.Cells.Item(89, 8).Value = 2319.8
$ws.Cells.Item(89, 9).Value = 2319.8
$ws.Cells.Item(89, 11).Value = 11599
$ws.Cells.Item(89, 13).Value = -5983

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 1850
$ws.Cells.Item(107, 9).Value = 1800
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 1800
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = 120
$ws.Cells.Item(107, 14).Value = -5840

$ws.Cells.Item(134, 8).Value = 3701.8333
$ws.Cells.Item(134, 9).Value = 4403.6665
$ws.Cells.Item(134, 11).Value = 13210.9995
$ws.Cells.Item(134, 13).Value = -10675.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 308.53333
$ws.Cells.Item(7, 9).Value = 326.76923
$ws.Cells.Item(7, 10).Value = 190
$ws.Cells.Item(7, 11).Value = 326.76923
$ws.Cells.Item(7, 12).Value = 190
$ws.Cells.Item(7, 13).Value = -213.76923
$ws.Cells.Item(7, 14).Value = -416

$ws.Cells.Item(16, 8).Value = 1118.8
$ws.Cells.Item(16, 9).Value = 645.25
$ws.Cells.Item(16, 10).Value = 3013
$ws.Cells.Item(16, 11).Value = 645.25
$ws.Cells.Item(16, 12).Value = 3013
$ws.Cells.Item(16, 13).Value = -358.25
$ws.Cells.Item(16, 14).Value = -3587

$ws.Cells.Item(31, 8).Value = 8109.579
$ws.Cells.Item(31, 9).Value = 1810.4286
$ws.Cells.Item(31, 10).Value = 25747.2
$ws.Cells.Item(31, 11).Value = 1810.4286
$ws.Cells.Item(31, 12).Value = 25747.2
$ws.Cells.Item(31, 13).Value = -1515.4286
$ws.Cells.Item(31, 14).Value = -26337.2

$ws.Cells.Item(34, 8).Value = 8109.579
$ws.Cells.Item(34, 9).Value = 1810.4286
$ws.Cells.Item(34, 10).Value = 25747.2
$ws.Cells.Item(34, 11).Value = 1810.4286
$ws.Cells.Item(34, 12).Value = 25747.2
$ws.Cells.Item(34, 13).Value = -1608.4286
$ws.Cells.Item(34, 14).Value = -26151.2

$ws.Cells.Item(107, 8).Value = 524.7273
$ws.Cells.Item(107, 9).Value = 470.33334
$ws.Cells.Item(107, 10).Value = 590
$ws.Cells.Item(107, 11).Value = 470.33334
$ws.Cells.Item(107, 12).Value = 590
$ws.Cells.Item(107, 13).Value = 1449.66666
$ws.Cells.Item(107, 14).Value = -4430

$ws.Cells.Item(113, 8).Value = 1118.8
$ws.Cells.Item(113, 9).Value = 645.25
$ws.Cells.Item(113, 10).Value = 3013
$ws.Cells.Item(113, 11).Value = 645.25
$ws.Cells.Item(113, 12).Value = 3013
$ws.Cells.Item(113, 13).Value = 1524.75
$ws.Cells.Item(113, 14).Value = -7353

$ws.Cells.Item(132, 8).Value = 2757.3157
$ws.Cells.Item(132, 9).Value = 2663.5
$ws.Cells.Item(132, 10).Value = 2918.1428
$ws.Cells.Item(132, 11).Value = 7990.5
$ws.Cells.Item(132, 12).Value = 8754.428400000001
$ws.Cells.Item(132, 13).Value = -5460.5
$ws.Cells.Item(132, 14).Value = -13814.4284

$ws.Cells.Item(134, 8).Value = 3928.2942
$ws.Cells.Item(134, 9).Value = 4384.7856
$ws.Cells.Item(134, 10).Value = 1798
$ws.Cells.Item(134, 11).Value = 13154.3568
$ws.Cells.Item(134, 12).Value = 5394
$ws.Cells.Item(134, 13).Value = -10619.3568
$ws.Cells.Item(134, 14).Value = -10464

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 70.46666999999999
$ws.Cells.Item(14, 9).Value = 70.46666999999999
$ws.Cells.Item(14, 11).Value = 211.40001
$ws.Cells.Item(14, 13).Value = -38.40000999999998

$ws.Cells.Item(107, 8).Value = 434.8
$ws.Cells.Item(107, 9).Value = 180
$ws.Cells.Item(107, 10).Value = 498.5
$ws.Cells.Item(107, 11).Value = 540
$ws.Cells.Item(107, 12).Value = 1495.5
$ws.Cells.Item(107, 13).Value = 1380
$ws.Cells.Item(107, 14).Value = -5335.5

$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 13).ClearContents()

$ws.Cells.Item(131, 8).Value = 36667544
$ws.Cells.Item(131, 10).Value = 47620090
$ws.Cells.Item(131, 12).Value = 142860270
$ws.Cells.Item(131, 14).Value = -142870350

$ws.Cells.Item(132, 8).Value = 1865.5883
$ws.Cells.Item(132, 9).Value = 1840.7273
$ws.Cells.Item(132, 10).Value = 1870.386
$ws.Cells.Item(132, 11).Value = 16566.5457
$ws.Cells.Item(132, 12).Value = 16833.474
$ws.Cells.Item(132, 13).Value = -14036.5457
$ws.Cells.Item(132, 14).Value = -21893.474

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 250001060
$ws.Cells.Item(113, 9).Value = 333334200
$ws.Cells.Item(113, 10).Value = 1600
$ws.Cells.Item(113, 11).Value = 333334200
$ws.Cells.Item(113, 12).Value = 1600
$ws.Cells.Item(113, 13).Value = -333332030
$ws.Cells.Item(113, 14).Value = -5940

$ws.Cells.Item(122, 8).Value = 2820248.2
$ws.Cells.Item(122, 9).Value = 5403217.5
$ws.Cells.Item(122, 10).Value = 2463.6365
$ws.Cells.Item(122, 11).Value = 16209652.5
$ws.Cells.Item(122, 12).Value = 7390.9095
$ws.Cells.Item(122, 13).Value = -16207202.5
$ws.Cells.Item(122, 14).Value = -12290.9095

$ws.Cells.Item(132, 8).Value = 4243.9165
$ws.Cells.Item(132, 9).Value = 3299
$ws.Cells.Item(132, 10).Value = 4329.8184
$ws.Cells.Item(132, 11).Value = 9897
$ws.Cells.Item(132, 12).Value = 12989.4552
$ws.Cells.Item(132, 13).Value = -7367
$ws.Cells.Item(132, 14).Value = -18049.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 27146856
$ws.Cells.Item(122, 9).Value = 35720284
$ws.Cells.Item(122, 11).Value = 107160852
$ws.Cells.Item(122, 13).Value = -107158402

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 2149.6667
$ws.Cells.Item(113, 9).Value = 1371.6666
$ws.Cells.Item(113, 10).Value = 3316.6667
$ws.Cells.Item(113, 11).Value = 4114.9998
$ws.Cells.Item(113, 12).Value = 9950.000100000001
$ws.Cells.Item(113, 13).Value = -1944.9998
$ws.Cells.Item(113, 14).Value = -14290.0001

$ws.Cells.Item(123, 8).Value = 29884.2
$ws.Cells.Item(123, 10).Value = 29884.2
$ws.Cells.Item(123, 12).Value = 29884.2
$ws.Cells.Item(123, 14).Value = -39684.2

$ws.Cells.Item(136, 8).Value = 1139.1034
$ws.Cells.Item(136, 9).Value = 776.41174
$ws.Cells.Item(136, 10).Value = 1652.9166
$ws.Cells.Item(136, 11).Value = 2329.23522
$ws.Cells.Item(136, 12).Value = 4958.7498
$ws.Cells.Item(136, 13).Value = 220.76478
$ws.Cells.Item(136, 14).Value = -10058.7498
